$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userCredentials")

$ws.Range("A7").Value = "FHDCMobileLogin"
$ws.Range("B7").Value = "testcaregiver2may@yopmail.com"
$ws.Range("C7").Value = "Qwerty@123"

$ws.Range("B10").Select()
